$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Special")
$ws.Activate()

# --- New weapon profile: Assault Bolter / Plasma (rows 78-79) ---
$ws.Range("A78").Value = "Assault Bolter"
$ws.Range("B78").Value = 3
$ws.Range("C78").Value = 18
$ws.Range("E78").Value = 4
$ws.Range("F78").Value = 10
$ws.Range("I78").Formula = "=(2/3)*`$B78*((13-E78)/12)"
$ws.Range("J78").Formula = "=(2/3)*`$B78*((13-F78)/12)"
$ws.Range("K78").Formula = "=I78+J78"

$ws.Range("A79").Value = "Plasma"
$ws.Range("B79").Value = 3
$ws.Range("C79").Value = 18
$ws.Range("E79").Value = 7
$ws.Range("F79").Value = 7
$ws.Range("I79").Formula = "=(2/3)*`$B79*((13-E79)/12)"
$ws.Range("J79").Formula = "=(2/3)*`$B79*((13-F79)/12)"
$ws.Range("K79").Formula = "=I79+J79"

# --- New weapon profiles: Heavy Onslaught / Plasma / Laser Destroyer (rows 81-83) ---
$ws.Range("A81").Value = "Heavy Onslaught"
$ws.Range("B81").Value = 3
$ws.Range("C81").Value = 30
$ws.Range("E81").Value = 6
$ws.Range("F81").Value = 9
$ws.Range("I81").Formula = "=(2/3)*`$B81*((13-E81)/12)"
$ws.Range("J81").Formula = "=(2/3)*`$B81*((13-F81)/12)"
$ws.Range("K81").Formula = "=I81+J81"

$ws.Range("A82").Value = "Plasma"
$ws.Range("B82").Value = 2
$ws.Range("C82").Value = 36
$ws.Range("E82").Value = 4
$ws.Range("F82").Value = 5
$ws.Range("I82").Formula = "=(2/3)*`$B82*((13-E82)/12)"
$ws.Range("J82").Formula = "=(2/3)*`$B82*((13-F82)/12)"
$ws.Range("K82").Formula = "=I82+J82"

$ws.Range("A83").Value = "Laser Destroyer"
$ws.Range("B83").Value = 3
$ws.Range("C83").Value = 72
$ws.Range("E83").Value = 11
$ws.Range("F83").Value = 5
$ws.Range("I83").Formula = "=(2/3)*`$B83*((13-E83)/12)"
$ws.Range("J83").Formula = "=(2/3)*`$B83*((13-F83)/12)"
$ws.Range("K83").Formula = "=I83+J83"

# --- Restore view state (frozen header row, scrolled down, active cell E83) ---
$ws.Range("A65").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("E83").Select()
